# DatoObjetivo.xlsx fix: move the "deportista" column (originally the last
# column, P) so it becomes the first column (A), shifting every other
# column one place to the right. This mirrors "cut column P, insert the
# cut cells before column A" in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cut the last used column (P, "deportista") and insert it before column A.
# This shifts old A:O -> B:P and drops the cut column's old content into
# the freshly inserted column A.
$ws.Columns.Item(16).Cut()
$ws.Columns.Item(1).Insert()

# The column that used to be P (now pushed one further right to Q) is left
# behind as a blank, formatted cell once the cut range is removed - give it
# an explicit (default) number format so it is persisted instead of being
# dropped as a fully empty cell.
$ws.Range("Q1").NumberFormat = "General"

# Leave column A selected, matching the resulting UI state after the insert.
$ws.Columns.Item(1).Select() | Out-Null
